# Generate Report for Handback
# Adds a new handback row (47177279-49fd-4440-bccf-7b00af8eac4b.md) to the
# Overview, zh-cn and de-de sheets of the handback-status workbook.

$wb = $excel.ActiveWorkbook

$fileId = "47177279-49fd-4440-bccf-7b00af8eac4b"
$fileName = "$fileId.md"
$pathAndName = "e2e\$fileId.md"
$ext = ".md"
$statusInSync = "Handed back: in sync with en-US"

$zhXlf = "$fileId.96e74346a401fcc5d27eb70e13f501e98a50543b.zh-cn.xlf"
$deXlf = "$fileId.96e74346a401fcc5d27eb70e13f501e98a50543b.de-de.xlf"

$zhHandoffDate = "2016-09-01 10:51:40"
$zhHandbackDate = "2016-09-01 10:51:57"
$deHandoffDate = "2016-09-01 10:51:44"
$deHandbackDate = "2016-09-01 10:52:14"

$mdHyperlinkBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0000000000000000000000000000000000000000/e2e/$fileName"
$zhHyperlinkBase = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0000000000000000000000000000000000000000/e2e/$fileName"
$deHyperlinkBase = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0000000000000000000000000000000000000000/e2e/$fileName"

# ---------------------------------------------------------------------------
# Overview sheet -> new row 4
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $fileName
$wsOverview.Range("B4").Value = $pathAndName
$wsOverview.Range("C4").Value = $ext
$wsOverview.Range("E4").Value = $statusInSync
$wsOverview.Range("F4").Value = $statusInSync
$wsOverview.Range("G4").Value = $deHandoffDate
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $mdHyperlinkBase, "", "", $pathAndName)

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# zh-cn sheet -> new row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $fileName
$wsZh.Range("B4").Value = $ext
$wsZh.Range("C4").Value = $statusInSync
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = $zhHandoffDate
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = $fileName
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").Value = $zhHandbackDate
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $mdHyperlinkBase, "", "", $fileName)
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), $zhHyperlinkBase, "", "", $fileName)

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P4"))

# ---------------------------------------------------------------------------
# de-de sheet -> new row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $fileName
$wsDe.Range("B4").Value = $ext
$wsDe.Range("C4").Value = $statusInSync
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = $deHandoffDate
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = $fileName
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").Value = $deHandbackDate
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $mdHyperlinkBase, "", "", $fileName)
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), $deHyperlinkBase, "", "", $fileName)

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P4"))
